$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 419, shifting rows 419:471 down to 420:472
$ws.Rows.Item(419).Insert()

# Populate the new row 419 with the inserted data
$ws.Cells.Item(419, 1).Value = 9
$ws.Cells.Item(419, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(419, 3).Value = "Metropolitana"
$ws.Cells.Item(419, 4).Value = 45142
$ws.Cells.Item(419, 5).Value = 13
$ws.Cells.Item(419, 6).Value = 300000001
$ws.Cells.Item(419, 7).Value = "Rabanito"
$ws.Cells.Item(419, 8).Value = "Sin especificar"
$ws.Cells.Item(419, 9).Value = "Primera"
$ws.Cells.Item(419, 10).Value = 7000
$ws.Cells.Item(419, 11).Value = 3000
$ws.Cells.Item(419, 12).Value = 3000
$ws.Cells.Item(419, 13).Value = 3000
$ws.Cells.Item(419, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(419, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(419, 16).Value = 30
$ws.Cells.Item(419, 17).Value = 100
$ws.Cells.Item(419, 18).Value = "Hortaliza"
